# Word COM-interop script applying the edits described by the commit diff.
# The underlying changes are all proofing-markup cleanups (removal of every
# <w:proofErr/> element left over from spell/grammar check) plus a handful of
# genuine textual corrections:
#   - "... sortir de la prison." -> "... sortir de la prison du château."
#   - "trois pièces, ou il devra" -> "trois pièces, où il devra"
#   - "interface graphique ou le joueur" -> "interface graphique où le joueur"
#   - "orienté objet" -> "orientée objet" (agreement)
#   - "conseillé pour utiliser" -> "conseillée pour utiliser" (agreement)
#   - "de fois ou elles" -> "de fois où elles"
#   - "qu'on à fini" -> "qu'on a fini"

$d = $word.ActiveDocument

# 1) "... sortir du château" (escape game intro paragraph) - only proofErr
#    cleanup around "game"/"games"; perform a no-op replace so the runs that
#    used to be split around the spell-checked words get merged together,
#    exactly like accepting/clearing the spelling flags would do in Word.
$d.Content.Find.Execute("est un mini escape game qui", $true, $false, $false, $false, $false, $true, 1, $false, "est un mini escape game qui", 2)
$d.Content.Find.Execute("Les escapes games sont des", $true, $false, $false, $false, $false, $true, 1, $false, "Les escapes games sont des", 2)

# 2) Insert "du château" after "réussir à sortir de la prison"
$d.Content.Find.Execute("réussir à sortir de la prison. Une fois sorti", $true, $false, $false, $false, $false, $true, 1, $false, "réussir à sortir de la prison du château. Une fois sorti", 2)

# 3) "ou" -> "où" : "un long couloir avec trois pièces, ou il devra"
$d.Content.Find.Execute("avec trois pièces, ou il devra", $true, $false, $false, $false, $false, $true, 1, $false, "avec trois pièces, où il devra", 2)

# 4) "ou" -> "où" : "interface graphique ou le joueur peut facilement interagir"
$d.Content.Find.Execute("interface graphique ou le joueur peut facilement interagir", $true, $false, $false, $false, $false, $true, 1, $false, "interface graphique où le joueur peut facilement interagir", 2)

# 5) "orienté" -> "orientée" and "conseillé" -> "conseillée" (feminine agreement)
$d.Content.Find.Execute("la programmation orienté objet", $true, $false, $false, $false, $false, $true, 1, $false, "la programmation orientée objet", 2)
$d.Content.Find.Execute("elle est conseillé pour utiliser le module Pyxel", $true, $false, $false, $false, $false, $true, 1, $false, "elle est conseillée pour utiliser le module Pyxel", 2)

# 6) Remaining proofErr-only cleanups (no text change) for "draw" occurrences
$d.Content.Find.Execute("la méthode update et la méthode draw. La méthode update met", $true, $false, $false, $false, $false, $true, 1, $false, "la méthode update et la méthode draw. La méthode update met", 2)
$d.Content.Find.Execute("ex : appuyer sur une touche). La méthode draw", $true, $false, $false, $false, $false, $true, 1, $false, "ex : appuyer sur une touche). La méthode draw", 2)
$d.Content.Find.Execute("La méthode draw :", $true, $false, $false, $false, $false, $true, 1, $false, "La méthode draw :", 2)

# 7) proofErr-only cleanup (no text change): "(le moment ou se trouve le joueur dans le jeu)"
$d.Content.Find.Execute("le moment ou se trouve le joueur dans le jeu", $true, $false, $false, $false, $false, $true, 1, $false, "le moment ou se trouve le joueur dans le jeu", 2)

# 8) proofErr-only cleanup (no text change): "en premier il est créer avec un tableau"
$d.Content.Find.Execute("en premier il est créer avec un tableau de tableau qui contient", $true, $false, $false, $false, $false, $true, 1, $false, "en premier il est créer avec un tableau de tableau qui contient", 2)

# 9) "ou" -> "où" : "le nombre de fois ou elles se trouvent"
$d.Content.Find.Execute("de fois ou elles se trouvent", $true, $false, $false, $false, $false, $true, 1, $false, "de fois où elles se trouvent", 2)

# 10) "à" -> "a" : "une fois qu'on à fini"
$d.Content.Find.Execute("une fois qu’on à fini", $true, $false, $false, $false, $false, $true, 1, $false, "une fois qu’on a fini", 2)
